# "Added support for longer quotes, fixed surplus numnber"
#
# The "surplus" surcharge multiplier (column K, "SS @ 1.0565" source rows)
# was set to 1.0565 on several line items by mistake - it should be 1
# (i.e. no surplus surcharge) for those rows. Correct each one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$surplusRows = @(16, 17, 20, 23, 25, 28, 31, 34, 35, 39)
foreach ($row in $surplusRows) {
    $ws.Cells.Item($row, 11).Value = 1
}

# Move/leave the cursor on K40 (the "Surplus" legend cell), matching where
# the author's selection ended up after making the edit.
$ws.Range("K40").Select() | Out-Null
